# Auto-generated Excel COM-interop script applying the Coeurl_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across
# several rows in the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1460
$ws.Range("I40").Value = 1196.5
$ws.Range("K40").Value = 1196.5
$ws.Range("M40").Value = -1021.5

# Row 51
$ws.Range("H51").Value = 2553.5356
$ws.Range("J51").Value = 3249.5715
$ws.Range("L51").Value = 3249.5715
$ws.Range("N51").Value = -4217.5715

# Row 68
$ws.Range("H68").Value = 24000
$ws.Range("I68").Value = 24000
$ws.Range("K68").Value = 24000
$ws.Range("M68").Value = -23251

# Row 71
$ws.Range("H71").Value = 24000
$ws.Range("I71").Value = 24000
$ws.Range("K71").Value = 72000
$ws.Range("M71").Value = -68256

# Row 86
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# Row 132
$ws.Range("H132").Value = 1646.9
$ws.Range("I132").Value = 1535.75
$ws.Range("J132").Value = 3203
$ws.Range("K132").Value = 4607.25
$ws.Range("L132").Value = 9609
$ws.Range("M132").Value = -2077.25
$ws.Range("N132").Value = -14669

# Row 138
$ws.Range("H138").Value = 6175364
$ws.Range("I138").Value = 1117.826
$ws.Range("K138").Value = 3353.478
$ws.Range("M138").Value = 1786.522

$ws = $wb.Worksheets.Item("ARM")
# Row 13
$ws.Range("H13").Value = 1946
$ws.Range("I13").Value = 1893
$ws.Range("J13").Value = 1999
$ws.Range("K13").Value = 1893
$ws.Range("L13").Value = 1999
$ws.Range("M13").Value = -1749
$ws.Range("N13").Value = -2287

# Row 102
$ws.Range("H102").Value = 2332.84
$ws.Range("I102").Value = 1963.9474
$ws.Range("K102").Value = 1963.9474
$ws.Range("M102").Value = -341.9474

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 2224.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2224.5
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2450.5

# Row 86
$ws.Range("H86").Value = 2314.838
$ws.Range("I86").Value = 2237.4517
$ws.Range("K86").Value = 2237.4517
$ws.Range("M86").Value = -1114.4517

# Row 89
$ws.Range("H89").Value = 2314.838
$ws.Range("I89").Value = 2237.4517
$ws.Range("K89").Value = 11187.2585
$ws.Range("M89").Value = -5571.2585

# Row 97
$ws.Range("H97").Value = 18000.5
$ws.Range("I97").Value = 18000.5
$ws.Range("K97").Value = 18000.5
$ws.Range("M97").Value = -17009.5

# Row 105
$ws.Range("H105").Value = 997.1111
$ws.Range("I105").Value = 998.4400000000001
$ws.Range("J105").Value = 980.5
$ws.Range("K105").Value = 998.4400000000001
$ws.Range("L105").Value = 980.5
$ws.Range("M105").Value = 748.5599999999999
$ws.Range("N105").Value = -4474.5

# Row 134
$ws.Range("H134").Value = 1509.8315
$ws.Range("I134").Value = 1496.076
$ws.Range("J134").Value = 1931.6666
$ws.Range("K134").Value = 4488.228
$ws.Range("L134").Value = 5794.9998
$ws.Range("M134").Value = -1953.228
$ws.Range("N134").Value = -10864.9998

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 188.7
$ws.Range("I22").Value = 188.7
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 188.7
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 62
$ws.Range("H62").Value = 11492.75
$ws.Range("J62").Value = 10077
$ws.Range("L62").Value = 10077
$ws.Range("N62").Value = -11325

# Row 65
$ws.Range("H65").Value = 11492.75
$ws.Range("J65").Value = 10077
$ws.Range("L65").Value = 50385
$ws.Range("N65").Value = -56625

# Row 97
$ws.Range("H97").Value = 38500
$ws.Range("J97").Value = 38500
$ws.Range("L97").Value = 38500
$ws.Range("N97").Value = -40482

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 24380.4
$ws.Range("I87").Value = 20634.166
$ws.Range("K87").Value = 61902.49800000001
$ws.Range("M87").Value = -60654.49800000001

# Row 90
$ws.Range("H90").Value = 24380.4
$ws.Range("I90").Value = 20634.166
$ws.Range("K90").Value = 185707.494
$ws.Range("M90").Value = -179467.494

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 4264.375
$ws.Range("J2").Value = 192.71428
$ws.Range("L2").Value = 192.71428
$ws.Range("N2").Value = -418.71428

# Row 97
$ws.Range("H97").Value = 412.5
$ws.Range("I97").Value = 517.6923
$ws.Range("J97").Value = 139
$ws.Range("K97").Value = 517.6923
$ws.Range("L97").Value = 139
$ws.Range("M97").Value = -21.69230000000005
$ws.Range("N97").Value = -1131

# Row 102
$ws.Range("H102").Value = 43481030
$ws.Range("I102").Value = 2738.1
$ws.Range("K102").Value = 2738.1
$ws.Range("M102").Value = -1116.1

# Row 113
$ws.Range("H113").Value = 2030.0454
$ws.Range("I113").Value = 2214.2222
$ws.Range("J113").Value = 1201.25
$ws.Range("K113").Value = 2214.2222
$ws.Range("L113").Value = 1201.25
$ws.Range("M113").Value = -44.22220000000016
$ws.Range("N113").Value = -5541.25

# Row 132
$ws.Range("H132").Value = 3173.4707
$ws.Range("I132").Value = 2465
$ws.Range("J132").Value = 5476
$ws.Range("K132").Value = 7395
$ws.Range("L132").Value = 16428
$ws.Range("M132").Value = -4865
$ws.Range("N132").Value = -21488

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2682.9333
$ws.Range("I61").Value = 2573.2964
$ws.Range("J61").Value = 3669.6667
$ws.Range("K61").Value = 2573.2964
$ws.Range("L61").Value = 3669.6667
$ws.Range("M61").Value = -2371.2964
$ws.Range("N61").Value = -4073.6667

# Row 113
$ws.Range("H113").Value = 2682.9333
$ws.Range("I113").Value = 2573.2964
$ws.Range("J113").Value = 3669.6667
$ws.Range("K113").Value = 2573.2964
$ws.Range("L113").Value = 3669.6667
$ws.Range("M113").Value = -403.2964000000002
$ws.Range("N113").Value = -8009.6667

# Row 132
$ws.Range("H132").Value = 5387.375
$ws.Range("I132").Value = 4871.2856
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 14613.8568
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -12083.8568
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 37899
$ws.Range("J93").Value = 37899
$ws.Range("L93").Value = 37899
$ws.Range("N93").Value = -42891

# Row 122
$ws.Range("H122").Value = 1704.1111
$ws.Range("I122").Value = 1538
$ws.Range("K122").Value = 4614
$ws.Range("M122").Value = -2164

# Row 126
$ws.Range("H126").Value = 13390.7
$ws.Range("J126").Value = 6999.5
$ws.Range("L126").Value = 20998.5
$ws.Range("N126").Value = -25938.5

# Row 130
$ws.Range("H130").Value = 61713
$ws.Range("J130").Value = 61713
$ws.Range("L130").Value = 61713
$ws.Range("N130").Value = -71753

# Row 138
$ws.Range("H138").Value = 55500
$ws.Range("J138").Value = 55500
$ws.Range("L138").Value = 55500
$ws.Range("N138").Value = -65780
